$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(56.426,56.426,3.03340581,0.04028454,2.26997876,129.56059548,6.449971017464042,6.449971017464042,0.3402677275178092,0.003814876311188826,0.318706954806613,31.02704578393793),
    @(86.163,86.163,1.99537222,0.03343293999999999,2.84798751,246.73747622,11.32257177639566,11.32257177639566,0.2629915256928592,0.004759148206144462,0.3532168311495181,51.72210462156272),
    @(31.54,63.035,2.729127,0.03505555,0.5526313900000001,17.74162784,4.284244945699855,8.544634912375471,0.3610348958957492,0.005039072190993699,0.1103469623427223,5.601350248477205),
    @(47.451,92.688,1.86378554,0.02807004,0.65785667,31.62002765,7.578194245862798,13.81243785287157,0.2769510662525184,0.004653765500294544,0.1130795881397919,9.24464776399464),
    @(18.105,72.288,2.42129739,0.02495465,0.11301532,2.11541271,3.242914043185772,12.93384476249581,0.4837121954432289,0.004824932357012868,0.03072817692207845,0.8867501776835333),
    @(26.451,96.639,1.78548637,0.02273616,0.1491138,4.053499850000001,5.020637987926017,14.16408487337958,0.256407048601429,0.004544597564297792,0.03660547286682422,1.774858587307004),
    @(12.243,73.29000000000001,2.42486441,0.0206359,0.04246956,0.54903085,2.657633160166018,15.91720004451968,0.5720805643859963,0.004600885386133449,0.01475908414405078,0.2935980370229359),
    @(18.781,94.43899999999999,1.84088512,0.01839505,0.05737918,1.12824401,4.246772194746744,15.75063050678802,0.3156188089121507,0.003265759633451838,0.01606636979426174,0.5714073840994911),
    @(9.012,71.884,2.507227299999999,0.01660387,0.01896527,0.18379979,2.235811900597824,17.84323190806321,0.661523204476586,0.003252455070725357,0.006913931122789022,0.1091010990103269),
    @(14.56,88.339,1.97984356,0.01576841,0.02877109,0.44888761,3.83863714554204,16.40900162515484,0.3711075514913437,0.002775071313082367,0.009467816879744156,0.2841402757676236),
    @(7.083,70.577,2.58359017,0.01505144,0.01093338,0.08521947000000001,1.940070842746438,19.28839322268999,0.7364057312791176,0.003110658857046996,0.004576088937680907,0.05857671725885571),
    @(11.58,80.089,2.19018567,0.01445044,0.01685115000000001,0.21177302,3.202163933007467,15.70362005025219,0.4227209638447388,0.00265218383327504,0.006173582420070003,0.1475696315894535)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        # Columns C..N are columns 3..14
        $colNum = $j + 3
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}

$wb.Save()
